# Add a "Predicted_revenue" column (R) computed as Predicted_ROI * production_budget + production_budget,
# and fix the previously mis-scaled Predicted_ROI (Q) values (they were off by a power-of-ten
# because of a bad scientific-notation export).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Predicted_ROI values in column Q (they were off by a power-of-ten scale)
$ws.Range("Q2").Value = 1.1913565769888199
$ws.Range("Q3").Value = 2.16002703607501
$ws.Range("Q4").Value = 1.3290001487029699
$ws.Range("Q5").Value = 3.1370468088940902
$ws.Range("Q6").Value = 1.9892420570115601
$ws.Range("Q7").Value = 2.0310178535422199
$ws.Range("Q8").Value = 1.2066929862648399
$ws.Range("Q9").Value = 1.85725487468071

# Re-format the Predicted_ROI column to 2 decimal places
$ws.Range("Q2:Q9").NumberFormat = "0.00"

# Add the new header (copy the header-row style from the neighboring cell)
$ws.Range("R1").Value = "Predicted_revenue"
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)

# Add the Predicted_revenue formula for each data row
$ws.Range("R2:R9").Formula = "=Q2*L2 + L2"
# Formula auto-inherits Q's number format from the adjacent cell; strip that back off
# so the new column keeps the default/general format, same as a freshly added column.
$ws.Range("R2:R9").ClearFormats()

# Widen column Q to fit the "Predicted_ROI" values/label
$ws.Columns("Q").ColumnWidth = 29.33
